$wb = $excel.ActiveWorkbook
$ws19 = $wb.Worksheets.Item("MAR-2022")
$ws19.Copy([System.Reflection.Missing]::Value, $ws19)
$ws20 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws20.Name = "TEST"
$ws20.Range("D6:D7").UnMerge()
$ws20.Range("D13:D14").UnMerge()
$ws20.Range("D20:D21").UnMerge()
$ws20.Range("D27:D28").UnMerge()
$ws20.Range("D3:D4").Merge()
$ws20.Range("D10:D11").Merge()
$ws20.Range("D15:D16").Merge()
$ws20.Range("D17:D18").Merge()
$ws20.Range("D24:D25").Merge()
Write-Host "done"
